$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.802.38"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "1.649.75"
$ws.Range("E4").Value = "  +0.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.74"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.19"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "1.879.40"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "1.656.62"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.18"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "26.817.17"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.50"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("E21").Value = "  +1.78%  "
$ws.Range("E22").Value = "  +15.54%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.51"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.77"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.13"
$ws.Range("E28").Value = "  +3.83%  "
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("E31").Value = "  +2.20%  "
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.00"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "1.287.33"
$ws.Range("E34").Value = "  +3.96%  "
$ws.Range("E35").Value = "  +3.06%  "
$ws.Range("E36").Value = "  +3.13%  "
$ws.Range("E37").Value = "  +2.89%  "
$ws.Range("E38").Value = "  +5.98%  "
$ws.Range("E39").Value = "  +4.33%  "
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("E43").Value = "  +2.58%  "
$ws.Range("D44").Value = "1.789.31"
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.03"
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.71"
$ws.Range("E46").Value = "  +8.80%  "
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.79"
$ws.Range("E49").Value = "  +3.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0970"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("E51").Value = "  +0.78%  "
